# Refresh the coinranking snapshot on Sheet1: new Price (column D) and
# Volume(1h) (column E) values, matching the "Updated cryptos list" run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.149.10"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "3.117.04"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("E9").Value = "  +0.36%  "

$ws.Range("E10").Value = "  +0.22%  "

$ws.Range("E11").Value = "  -0.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000249"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("E14").Value = "  -1.56%  "

$ws.Range("D15").Value = "3.633.72"
$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("D16").Value = "67.124.51"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("E17").Value = "  -0.86%  "

$ws.Range("D18").Value = "3.118.24"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("E19").Value = "  +1.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "492.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.707"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("E22").Value = "  +4.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "84.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.43%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.82%  "

$ws.Range("E29").Value = "  -1.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.74%  "

$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").Value = "0.0₃0948"
$ws.Range("E33").Value = "  -5.64%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.976"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("E38").Value = "  -2.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.310"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.94%  "

$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "386.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.96%  "

$ws.Range("D43").Value = "2.818.15"
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("E44").Value = "  -6.80%  "

$ws.Range("E45").Value = "  -2.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.21%  "

$ws.Range("E49").Value = "  -0.93%  "

$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("E51").Value = "  -1.46%  "
